$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of data rows currently on the sheet (before inserting the new column)
$rowCount = $ws.UsedRange.Rows.Count

# Insert a brand-new column before column A; this pushes the existing
# question/answer columns (A, B) over to (B, C).
$ws.Columns("A").Insert()

# Fill the freshly inserted column A with a 0-based row index.
for ($i = 0; $i -lt $rowCount; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $i
}

# Move the active selection, matching the saved workbook view state.
[void]$ws.Range("E17").Select()
